# Updated symbol list on Mon Jan  9 09:14:09 UTC 2023 with GitHub Actions
#
# Applies refreshed crypto price / volume / hour data to the worksheet.
# All data cells in this sheet are stored as text (inline strings), including
# cells that look numeric (e.g. "277.45", "0.06282") or percentages
# (e.g. "6.26%"). We must preserve that text representation exactly -
# assigning a plain numeric-looking string to .Value would make Excel
# coerce it into a real number (losing trailing zeros / exact formatting).
# Prefixing the string with a leading apostrophe forces Excel to keep it
# as text, matching the original "t=\"inlineStr\"" cell type/content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $ws.Range($cell).Value = "'" + $value
}

# Rows 2 through 51 all have their "Hora" (column G) value updated from 8 to 9.
for ($row = 2; $row -le 51; $row++) {
    Set-TextValue "G$row" "9"
}

# Rows with updated Price (column D) and/or Volume(1h) (column E) values.
$updates = @(
    @{ Row = 2; D = "277.45"; E = "6.26%" },
    @{ Row = 3; D = "27.32"; E = "1.66%" },
    @{ Row = 4; D = "4.784"; E = "1.55%" },
    @{ Row = 5; D = "0.06282"; E = "1.10%" },
    @{ Row = 6; D = "6.930"; E = "3.02%" },
    @{ Row = 7; D = "0.8786"; E = "3.43%" },
    @{ Row = 8; D = "0.9449"; E = "3.45%" },
    @{ Row = 9; D = "0.1455"; E = "3.95%" },
    @{ Row = 10; D = "0.05165"; E = "4.45%" },
    @{ Row = 11; D = "0.07287"; E = "" },
    @{ Row = 12; D = "0.03098"; E = "-0.63%" },
    @{ Row = 13; D = "0.09067"; E = "0.18%" },
    @{ Row = 14; D = "0.001552"; E = "1.51%" },
    @{ Row = 15; D = "0.0006283"; E = "1.46%" },
    @{ Row = 16; D = "0.005839"; E = "-2.37%" },
    @{ Row = 17; D = "3.447"; E = "-0.07%" },
    @{ Row = 18; D = "3.271"; E = "3.14%" },
    @{ Row = 19; D = ""; E = "5.31%" },
    @{ Row = 21; D = "0.1313"; E = "" },
    @{ Row = 22; D = "3.854"; E = "-5.94%" },
    @{ Row = 23; D = "0.04329"; E = "2.41%" },
    @{ Row = 24; D = "0.001183"; E = "0.14%" },
    @{ Row = 25; D = "0.004274"; E = "5.06%" },
    @{ Row = 26; D = "0.0001202"; E = "" },
    @{ Row = 27; D = ""; E = "3.19%" },
    @{ Row = 40; D = "0.04066"; E = "3.12%" },
    @{ Row = 41; D = "0.006419"; E = "55.19%" },
    @{ Row = 42; D = "0.1155"; E = "3.83%" },
    @{ Row = 43; D = "0.002177"; E = "1.61%" },
    @{ Row = 44; D = "0.01186"; E = "-9.91%" },
    @{ Row = 45; D = "0.00005144"; E = "-0.39%" },
    @{ Row = 46; D = ""; E = "0.04%" },
    @{ Row = 47; D = ""; E = "858.56%" },
    @{ Row = 48; D = "0.02253"; E = "-33.79%" },
    @{ Row = 49; D = "0.00002103"; E = "0.04%" },
    @{ Row = 50; D = ""; E = "0.04%" }
)

foreach ($u in $updates) {
    if ($u.D -ne "") {
        Set-TextValue "D$($u.Row)" $u.D
    }
    if ($u.E -ne "") {
        Set-TextValue "E$($u.Row)" $u.E
    }
}
